$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.793.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.620.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0754"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.808.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0155"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.799"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.765.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0110"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.29%  "
